$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Fix the duplicated "d'accueil d'accueil" -> "d'accueil" and drop
#    the new cursor position bookmark ("_GoBack") right after the
#    fixed run (this is where Word leaves the _GoBack bookmark after
#    the last edit made by the author).
# ------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute(
    "Voici ci-dessous une image de la page d" + [char]0x2019 + "accueil d" + [char]0x2019 + "accueil",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Voici ci-dessous une image de la page d" + [char]0x2019 + "accueil",
    2)
Write-Output "Fix duplicated text: $found"

$fixEndPos = $rng.End
$goBackRange = $d.Range($fixEndPos, $fixEndPos)
$d.Bookmarks.Add("_GoBack", $goBackRange)

# ------------------------------------------------------------------
# 2) Merge the two runs that make up the "Comme nous pouvons..."
#    paragraph into a single run, while leaving the trailing "."
#    run (which belongs to a distinct sentence) untouched. A
#    temporary bookmark is used purely to stop that trailing run
#    from being folded into the merge, then it is removed again.
# ------------------------------------------------------------------
$tailRng = $d.Content
$tailFound = $tailRng.Find.Execute("Nous en reparlerons plus tard dans le documents.")
Write-Output "Locate trailing sentence: $tailFound"
$periodStart = $tailRng.End - 1
$periodRange = $d.Range($periodStart, $tailRng.End)
$d.Bookmarks.Add("zzTempProtect", $periodRange)

$mergeRng = $d.Content
$mergeFound = $mergeRng.Find.Execute(
    "plateforme web. Nous avons aussi",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "plateforme web. Nous avons aussi",
    2)
Write-Output "Merge runs: $mergeFound"

$tempBm = $d.Bookmarks.Item("zzTempProtect")
$tempBm.Delete()

# ------------------------------------------------------------------
# 3) Word only ever keeps a single "_GoBack" bookmark (it tracks the
#    last edit position). Having re-added it in step 1 at the new
#    edit location automatically relocated it away from its former
#    spot on the empty paragraph just before "Pour finir, le
#    bouton...", which now collapses back to a plain empty paragraph.
# ------------------------------------------------------------------
Write-Output "Bookmarks remaining: $($d.Bookmarks.Count)"

Write-Output "Done"
